$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 5
$ws.Range("D5").Value = 2
$ws.Range("F5").Value = -3
$ws.Range("H5").Value = 46

# Update the active selection to match the saved view state
$ws.Range("H8").Select()
